$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column O (currently "Extracted Objects" ...),
# which shifts old O:U (Extracted Objects .. Result String) to Q:W.
$ws.Range("O1:P1").EntireColumn.Insert()

# Rename the existing (now still M/N) predicate-with-parents/related headers.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Set headers for the two newly inserted columns.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Fill in the new column values (Correct Pred Predicates Parents / Related)
# for each data row.
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 2

$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4

$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 3

$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 4

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
